$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added for "Macroferia Regional de Talca" (Brócoli),
# inserted at row 49 -- pushing every subsequent record down by one row.
$ws.Rows(49).Insert()

$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44435
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = 100112023
$ws.Range("G49").Value = "Brócoli"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Segunda"
$ws.Range("J49").Value = 8000
$ws.Range("K49").Value = 450
$ws.Range("L49").Value = 5000
$ws.Range("M49").Value = 2175
$ws.Range("N49").Value = "$/unidad"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 2175
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
